$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column E (Runmode) for rows 8 through 21 from "Yes" to "No"
for ($r = 8; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "No"
}

# Update the selection shown in the sheet view to match the new active range
$ws.Range("E8:E21").Select()
